$d = $word.ActiveDocument

# Namespace prefix used for InsertXML single-part WordprocessingML packages.
$pkgHeader = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
$pkgFooter = '</w:document></pkg:xmlData></pkg:part></pkg:package>'

# ---------------------------------------------------------------------------
# 1) Split the "Socks in the Dark" intro sentence so that "following:" is its
#    own run wrapped in proofErr gramStart/gramEnd markers.
# ---------------------------------------------------------------------------
function Find-ParagraphIndex($doc, $matchText) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs($i).Range.Text -like $matchText) {
            return $i
        }
    }
    return -1
}

$idx = Find-ParagraphIndex $d "*guarantee getting the following:*"
$p = $d.Paragraphs($idx)
$r = $p.Range
$r2 = $d.Range($r.Start, $r.End - 1)

$body = '<w:body><w:p>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Calibri Bold Italic" w:hAnsi="Calibri Bold Italic" w:cs="Calibri Bold Italic"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr>' +
    '<w:t xml:space="preserve">There are 20 socks in a drawer: 5 pairs of black socks, 3 pairs of brown and 2 pairs of white. You select the socks in the dark and can check them only after a selection has been made. What is the smallest number of socks you need to select to guarantee getting the </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Calibri Bold Italic" w:hAnsi="Calibri Bold Italic" w:cs="Calibri Bold Italic"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>following:</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '</w:p></w:body>'

[void]$r2.InsertXML($pkgHeader + $body + $pkgFooter)

# ---------------------------------------------------------------------------
# 2) Add two new paragraphs after "...4+2=18" with the follow-up question and
#    its answer (the actual new content from the commit). The hidden
#    "_GoBack" bookmark that used to sit at the end of the "4+2=18"
#    paragraph moves along to the end of the new last paragraph.
# ---------------------------------------------------------------------------
$idx18 = Find-ParagraphIndex $d "*4+2=18*"
$p18 = $d.Paragraphs($idx18)
$r18 = $p18.Range

$body2 = '<w:body>' +
    '<w:p><w:pPr><w:widowControl w:val="0"/><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="240"/><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/></w:rPr></w:pPr>' +
        '<w:r><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/></w:rPr><w:t>20-</w:t></w:r>' +
        '<w:r><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/></w:rPr><w:t>4+2=18</w:t></w:r>' +
    '</w:p>' +
    '<w:p><w:pPr><w:widowControl w:val="0"/><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="240"/><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/></w:rPr></w:pPr>' +
        '<w:r><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/></w:rPr><w:t>Using the same thinking but changing the first selection to B what will be the answer?</w:t></w:r>' +
    '</w:p>' +
    '<w:p><w:pPr><w:widowControl w:val="0"/><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="240"/><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/></w:rPr></w:pPr>' +
        '<w:r><w:rPr><w:rFonts w:ascii="Times" w:hAnsi="Times" w:cs="Times"/></w:rPr><w:t>20-6+2=16</w:t></w:r>' +
        '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
    '</w:p>' +
    '</w:body>'

[void]$r18.InsertXML($pkgHeader + $body2 + $pkgFooter)

# ---------------------------------------------------------------------------
# 3) Split each of the "a)/b)/c) What if the girl counts..." paragraphs so
#    the numeric answer is wrapped by proofErr gramStart/gramEnd, and the
#    trailing space (for a/b) becomes its own run.
# ---------------------------------------------------------------------------
$idxA = Find-ParagraphIndex $d "*girl counts from 1 to 10 *"
$pA2 = $d.Paragraphs($idxA)
$rA = $pA2.Range
$rA2 = $d.Range($rA.Start, $rA.End - 1)
$bodyA = '<w:body><w:p>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Calibri Bold Italic" w:hAnsi="Calibri Bold Italic" w:cs="Calibri Bold Italic"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve">a) What if the girl counts from 1 to </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Calibri Bold Italic" w:hAnsi="Calibri Bold Italic" w:cs="Calibri Bold Italic"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>10</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Calibri Bold Italic" w:hAnsi="Calibri Bold Italic" w:cs="Calibri Bold Italic"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
    '</w:p></w:body>'
[void]$rA2.InsertXML($pkgHeader + $bodyA + $pkgFooter)

$idxB = Find-ParagraphIndex $d "*girl counts from 1 to 100 *"
$pB2 = $d.Paragraphs($idxB)
$rB = $pB2.Range
$rB2 = $d.Range($rB.Start, $rB.End - 1)
$bodyB = '<w:body><w:p>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Calibri Bold Italic" w:hAnsi="Calibri Bold Italic" w:cs="Calibri Bold Italic"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve">b) What if the girl counts from 1 to </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Calibri Bold Italic" w:hAnsi="Calibri Bold Italic" w:cs="Calibri Bold Italic"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>100</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Calibri Bold Italic" w:hAnsi="Calibri Bold Italic" w:cs="Calibri Bold Italic"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
    '</w:p></w:body>'
[void]$rB2.InsertXML($pkgHeader + $bodyB + $pkgFooter)

$idxC = Find-ParagraphIndex $d "*girl counts from 1 to 1000*"
$pC2 = $d.Paragraphs($idxC)
$rC = $pC2.Range
$rC2 = $d.Range($rC.Start, $rC.End - 1)
$bodyC = '<w:body><w:p>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Calibri Bold Italic" w:hAnsi="Calibri Bold Italic" w:cs="Calibri Bold Italic"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve">c) What if the girl counts from 1 to </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Calibri Bold Italic" w:hAnsi="Calibri Bold Italic" w:cs="Calibri Bold Italic"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>1000</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '</w:p></w:body>'
[void]$rC2.InsertXML($pkgHeader + $bodyC + $pkgFooter)

Write-Host "Done."
